# Re-run with updated workers/jobs balance figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated NetMigration values (column C) ---
$ws.Range("C2").Value  = 50000
$ws.Range("C6").Value  = 300000
$ws.Range("C7").Value  = 325000
$ws.Range("C8").Value  = 15000
$ws.Range("C9").Value  = 30000
$ws.Range("C10").Value = 25000
$ws.Range("C11").Value = 25000
$ws.Range("C12").Value = 25000
$ws.Range("C13").Value = 15000
$ws.Range("C19").Value = 15000
$ws.Range("C21").Value = 15000
$ws.Range("C22").Value = 10000
$ws.Range("C24").Value = 5000

# --- Column widths (nudged slightly wider) ---
$ws.Columns("A").ColumnWidth = 12.3
$ws.Columns("C").ColumnWidth = 13.3

# --- View / selection: scroll back to top-left, select F8 instead of F19 ---
$ws.Range("F8").Select() | Out-Null
